$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entry for 2018-05-09
$ws.Range("A5").Value = 43229
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats -> reuse the date style from A1

$ws.Range("B5").Value = "Commencer à faire le MCD et le MLD. Je les ai montrés à M. Chavey et on en a parlé puis amené quelques modifications."
$ws.Range("C5").Value = "4 heures "

$ws.Range("C6").Select()
